# Add data for 2022-04-10
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-04-02"

# Update the label text for the April row
$ws.Range("A5").Value = "April (through 04-02)"

# Update row 5 (April) values
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 5
$ws.Range("I5").Value = 6

# Update row 6 (Total) values
$ws.Range("C6").Value = 130
$ws.Range("D6").Value = 191
$ws.Range("E6").Value = 202
$ws.Range("F6").Value = 113
$ws.Range("G6").Value = 202
$ws.Range("H6").Value = 428
$ws.Range("I6").Value = 439
